# Reorders the "Out of PO" player table and adds a new row for
# Robert Williams III, pushing Bradley Beal down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Payton Pritchard", "PG", "Boston Celtics"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Marcus Smart", "PG,SG", "Memphis Grizzlies"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
